# Project Plan.docx edit:
#   - Title placeholder text "<Project Name>test" -> "Project4"
#   - The old "Top" bookmark that wrapped "<Project Name>" is removed along
#     with the text it marked (the replacement text is no longer bookmarked).

$d = $word.ActiveDocument

# Remove the "Top" bookmark first (before editing the range) so Word does
# not just shrink/relocate it onto the surviving text.
if ($d.Bookmarks.Exists("Top")) {
    $d.Bookmarks("Top").Delete()
}

# Locate the Title-styled paragraph (first one in this document) and
# replace its whole text with "Project4", keeping the paragraph mark and
# its run formatting (yellow highlight / en-US language) intact.
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Title") {
        $titlePara = $p
        break
    }
}

if ($titlePara -ne $null) {
    $r = $titlePara.Range
    $r.MoveEnd(1, -1)
    $r.Text = "Project4"
}
